# "Fruta / hortaliza, semanal"
#
# A new weekly price-report row for Cilantro at "Terminal La Palmera de
# La Serena" is inserted as row 98 (pushing the previously existing rows
# 98-130 down to 99-131). The new row carries a new reporting date and
# new min/max/weighted-avg/per-kg price figures; every other attribute
# (market, region, product, unit, origin, classification, etc.) is the
# same boilerplate used by every other row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98; this shifts rows 98:130 down to
# 99:131 (carrying their formatting/styles with them, e.g. the date
# number format on column D), exactly like pressing Excel's
# "Insert Sheet Rows" on a selected row.
$ws.Rows(98).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A98").Value = 8
$ws.Range("B98").Value = "Terminal La Palmera de La Serena"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = 44642
$ws.Range("E98").Value = 4
$ws.Range("F98").Value = 100112040
$ws.Range("G98").Value = "Cilantro"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2300
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = 2400
$ws.Range("N98").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O98").Value = "Provincia del Elquí"
$ws.Range("P98").Value = 1600
$ws.Range("Q98").Value = 1.5
$ws.Range("R98").Value = "Hortaliza"
